# Fix Dates in Seminarplan
#
# 1) The registration/abstract deadline date "31.9.2023" is not a valid date
#    (September only has 30 days) -> fix to "30.9.2023".
# 2) The session ("Sitzung") numbering skips from "Sitzung 7" straight to
#    "Sitzung 9" (there never was a "Sitzung 8"). Renumber sessions 9-13
#    down by one so the sequence is contiguous: 9->8, 10->9, 11->10,
#    12->11, 13->12. Processed in ascending order so that a freshly
#    produced number (e.g. the new "Sitzung 9") is never re-matched by a
#    later replacement meant for the old "Sitzung 10", etc.

$d = $word.ActiveDocument

# --- 1) Fix the invalid deadline date 31.9.2023 -> 30.9.2023 ---------------
$rng = $d.Content
$rng.Find.Execute("31.9.2023", $true, $false, $false, $false, $false, $true, 1, $false, "30.9.2023", 2) | Out-Null

# --- 2) Renumber the sessions, lowest number first --------------------------
$rng = $d.Content
$rng.Find.Execute("Sitzung 9", $true, $false, $false, $false, $false, $true, 1, $false, "Sitzung 8", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Sitzung 10", $true, $false, $false, $false, $false, $true, 1, $false, "Sitzung 9", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Sitzung 11", $true, $false, $false, $false, $false, $true, 1, $false, "Sitzung 10", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Sitzung 12", $true, $false, $false, $false, $false, $true, 1, $false, "Sitzung 11", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Sitzung 13", $true, $false, $false, $false, $false, $true, 1, $false, "Sitzung 12", 2) | Out-Null
